# Auto-generated Excel COM-interop script to apply the Famfrit_Profits market-data refresh diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across all 8 sheets

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 448.9
$ws.Range("I2").Value = 440
$ws.Range("K2").Value = 440
$ws.Range("M2").Value = -327
$ws.Range("H34").Value = 3598.6667
$ws.Range("I34").Value = 3598.6667
$ws.Range("K34").Value = 3598.6667
$ws.Range("M34").Value = -3395.6667
$ws.Range("H36").Value = 3598.6667
$ws.Range("I36").Value = 3598.6667
$ws.Range("K36").Value = 3598.6667
$ws.Range("M36").Value = -2883.6667
$ws.Range("H38").Value = 4308.727
$ws.Range("I38").Value = 1056.7142
$ws.Range("K38").Value = 3170.1426
$ws.Range("M38").Value = -2798.1426
$ws.Range("H39").Value = 657.3333
$ws.Range("I39").Value = 309.77777
$ws.Range("J39").Value = 1700
$ws.Range("K39").Value = 929.33331
$ws.Range("L39").Value = 5100
$ws.Range("M39").Value = -633.33331
$ws.Range("N39").Value = -5692
$ws.Range("H53").Value = 504
$ws.Range("I53").Value = 443.14285
$ws.Range("J53").Value = 532.4
$ws.Range("K53").Value = 443.14285
$ws.Range("L53").Value = 532.4
$ws.Range("M53").Value = 193.85715
$ws.Range("N53").Value = -1806.4
$ws.Range("H80").Value = 2022
$ws.Range("J80").Value = 2481.6667
$ws.Range("L80").Value = 7445.000100000001
$ws.Range("N80").Value = -9441.000100000001
$ws.Range("H83").Value = 2022
$ws.Range("J83").Value = 2481.6667
$ws.Range("L83").Value = 22335.0003
$ws.Range("N83").Value = -32319.0003
$ws.Range("H94").Value = 2933.3333
$ws.Range("I94").Value = 2933.3333
$ws.Range("K94").Value = 2933.3333
$ws.Range("M94").Value = -2482.3333
$ws.Range("H115").Value = 909.6
$ws.Range("I115").Value = 909.6
$ws.Range("K115").Value = 2728.8
$ws.Range("M115").Value = -1161.8
$ws.Range("H129").Value = 1366.3334
$ws.Range("J129").Value = 1700
$ws.Range("L129").Value = 5100
$ws.Range("N129").Value = -15100
$ws.Range("H133").Value = 123520
$ws.Range("J133").Value = 123520
$ws.Range("L133").Value = 123520
$ws.Range("N133").Value = -133640

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 317.7143
$ws.Range("I4").Value = 329
$ws.Range("K4").Value = 329
$ws.Range("M4").Value = -213
$ws.Range("H17").Value = 910
$ws.Range("I17").Value = 1520
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 1520
$ws.Range("L17").Value = 300
$ws.Range("M17").Value = -1347
$ws.Range("N17").Value = -646
$ws.Range("H44").Value = 48711.25
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 48711.25
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 48711.25
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -49687.25
$ws.Range("H63").Value = 4386.875
$ws.Range("J63").Value = 8573
$ws.Range("L63").Value = 8573
$ws.Range("N63").Value = -9945
$ws.Range("H66").Value = 4386.875
$ws.Range("J66").Value = 8573
$ws.Range("L66").Value = 42865
$ws.Range("N66").Value = -49729
$ws.Range("H102").Value = 4722.7144
$ws.Range("I102").Value = 4009.8333
$ws.Range("K102").Value = 4009.8333
$ws.Range("M102").Value = -2387.8333
$ws.Range("H110").Value = 34534.184
$ws.Range("I110").Value = 37637.6
$ws.Range("K110").Value = 37637.6
$ws.Range("M110").Value = -35592.6
$ws.Range("H122").Value = 2820.76
$ws.Range("I122").Value = 2105.2144
$ws.Range("J122").Value = 3731.4546
$ws.Range("K122").Value = 6315.6432
$ws.Range("L122").Value = 11194.3638
$ws.Range("M122").Value = -3865.6432
$ws.Range("N122").Value = -16094.3638
$ws.Range("H132").Value = 55580852
$ws.Range("I132").Value = 13963
$ws.Range("J132").Value = 111147740
$ws.Range("K132").Value = 41889
$ws.Range("L132").Value = 333443220
$ws.Range("M132").Value = -39359
$ws.Range("N132").Value = -333448280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2256.2856
$ws.Range("I99").Value = 2282.3333
$ws.Range("J99").Value = 2100
$ws.Range("K99").Value = 2282.3333
$ws.Range("L99").Value = 2100
$ws.Range("M99").Value = -784.3332999999998
$ws.Range("N99").Value = -5096
$ws.Range("H132").Value = 75000
$ws.Range("J132").Value = 75000
$ws.Range("L132").Value = 75000
$ws.Range("N132").Value = -85120

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2246.4517
$ws.Range("I31").Value = 1417
$ws.Range("K31").Value = 1417
$ws.Range("M31").Value = -1122
$ws.Range("H34").Value = 2246.4517
$ws.Range("I34").Value = 1417
$ws.Range("K34").Value = 1417
$ws.Range("M34").Value = -1215
$ws.Range("H39").Value = 14011.4
$ws.Range("I39").Value = 7500
$ws.Range("K39").Value = 7500
$ws.Range("M39").Value = -7109
$ws.Range("H49").Value = 14011.4
$ws.Range("I49").Value = 7500
$ws.Range("K49").Value = 7500
$ws.Range("M49").Value = -7318
$ws.Range("H141").Value = 196775.33
$ws.Range("I141").Value = 30000
$ws.Range("K141").Value = 30000
$ws.Range("M141").Value = -24820

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 626.2857
$ws.Range("J107").Value = 646.3333
$ws.Range("L107").Value = 1938.9999
$ws.Range("N107").Value = -5778.9999
$ws.Range("H131").Value = 3038.0454
$ws.Range("J131").Value = 3499.7646
$ws.Range("L131").Value = 10499.2938
$ws.Range("N131").Value = -20579.2938
$ws.Range("H132").Value = 6156.091
$ws.Range("I132").Value = 2023
$ws.Range("J132").Value = 8517.857
$ws.Range("K132").Value = 18207
$ws.Range("L132").Value = 76660.713
$ws.Range("M132").Value = -15677
$ws.Range("N132").Value = -81720.713

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1016226.8
$ws.Range("J11").Value = 12916.667
$ws.Range("L11").Value = 12916.667
$ws.Range("N11").Value = -13194.667
$ws.Range("H43").Value = 4004.7273
$ws.Range("I43").Value = 1348.1111
$ws.Range("J43").Value = 15959.5
$ws.Range("K43").Value = 1348.1111
$ws.Range("L43").Value = 15959.5
$ws.Range("M43").Value = -1197.1111
$ws.Range("N43").Value = -16261.5
$ws.Range("H109").Value = 38003.8
$ws.Range("J109").Value = 36006.332
$ws.Range("L109").Value = 36006.332
$ws.Range("N109").Value = -38086.332
$ws.Range("H122").Value = 1900.2727
$ws.Range("I122").Value = 1535.3077
$ws.Range("J122").Value = 2427.4443
$ws.Range("K122").Value = 4605.9231
$ws.Range("L122").Value = 7282.3329
$ws.Range("M122").Value = -2155.9231
$ws.Range("N122").Value = -12182.3329
$ws.Range("H126").Value = 4545.154
$ws.Range("I126").Value = 4014.5
$ws.Range("K126").Value = 12043.5
$ws.Range("M126").Value = -9573.5
$ws.Range("H132").Value = 8982.200000000001
$ws.Range("I132").Value = 8475.817999999999
$ws.Range("J132").Value = 10374.75
$ws.Range("K132").Value = 25427.454
$ws.Range("L132").Value = 31124.25
$ws.Range("M132").Value = -22897.454
$ws.Range("N132").Value = -36184.25
$ws.Range("H133").Value = 104995
$ws.Range("J133").Value = 104995
$ws.Range("L133").Value = 104995
$ws.Range("N133").Value = -115115

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 595
$ws.Range("I46").Value = 595
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 595
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -407
$ws.Range("N46").ClearContents()
$ws.Range("H61").Value = 2563.0967
$ws.Range("J61").Value = 8978.666999999999
$ws.Range("L61").Value = 8978.666999999999
$ws.Range("N61").Value = -9382.666999999999
$ws.Range("H113").Value = 2563.0967
$ws.Range("J113").Value = 8978.666999999999
$ws.Range("L113").Value = 8978.666999999999
$ws.Range("N113").Value = -13318.667
$ws.Range("H122").Value = 4196.8335
$ws.Range("I122").Value = 3961.5833
$ws.Range("J122").Value = 4667.3335
$ws.Range("K122").Value = 11884.7499
$ws.Range("L122").Value = 14002.0005
$ws.Range("M122").Value = -9434.749899999999
$ws.Range("N122").Value = -18902.0005
$ws.Range("H133").Value = 55772.25
$ws.Range("J133").Value = 55772.25
$ws.Range("L133").Value = 55772.25
$ws.Range("N133").Value = -60832.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H107").Value = 1961
$ws.Range("I107").Value = 1559.8
$ws.Range("J107").Value = 2462.5
$ws.Range("K107").Value = 4679.4
$ws.Range("L107").Value = 7387.5
$ws.Range("M107").Value = -2759.4
$ws.Range("N107").Value = -11227.5
$ws.Range("H133").Value = 56211
$ws.Range("J133").Value = 56211
$ws.Range("L133").Value = 56211
$ws.Range("N133").Value = -66331
